# "Fruta / hortaliza, semanal" weekly update:
# A new weekly record is inserted as row 97 (Terminal Hortofrutícola Agro
# Chillán - Choclo, "Choclero" "Primera", 2021-12-29), which pushes every
# existing data row from 97..125 down to 98..126. The sheet's dimension
# grows from A1:R125 to A1:R126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 97, shifting rows 97-125
# (and everything below) down by one.
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly record.
$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(97, 3).Value = 'Ñuble'
$ws.Cells.Item(97, 4).Value = 44559
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = 100112024
$ws.Cells.Item(97, 7).Value = 'Choclo'
$ws.Cells.Item(97, 8).Value = 'Choclero'
$ws.Cells.Item(97, 9).Value = 'Primera'
$ws.Cells.Item(97, 10).Value = 14000
$ws.Cells.Item(97, 11).Value = 300
$ws.Cells.Item(97, 12).Value = 350
$ws.Cells.Item(97, 13).Value = 325
$ws.Cells.Item(97, 14).Value = '$/unidad'
$ws.Cells.Item(97, 15).Value = 'Región del Maule'
$ws.Cells.Item(97, 16).Value = 325
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = 'Hortaliza'

# Match the date number format used by the other rows in column D.
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
